$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item(1)
$wsData = $wb.Worksheets.Item(2)

# --- Update the metabolite identifiers / human-readable names on the
# "openbis-data" sheet (SE-341 metabolomics example data). ---
$wsData.Range("A2").Value = "CHEBI:15521"
$wsData.Range("A3").Value = "CHEBI:18311"
$wsData.Range("B2").Value = "phosphate1"
$wsData.Range("B3").Value = "phosphate2"

# --- Make "openbis-data" the active sheet/tab, with B8 selected there. ---
$wsData.Activate()
$wsData.Range("B8").Select()

# --- Leave the other sheet's own selection as previously recorded (C8). ---
$wsMeta.Range("C8").Select()
$wsData.Activate()
